# Auto-generated: applies scheduled-runner market-price refresh to Hyperion_Profits workbook
# Updates currentAveragePrice / NQ / HQ columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# to match the latest Universalis market snapshot used by the crafting-profit calculator.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4: Root Rush / Growth Formula Alpha
$ws.Range("H4").Value = 336.75
$ws.Range("I4").Value = 32.666668
$ws.Range("K4").Value = 32.666668
$ws.Range("M4").Value = 81.333332

# ALC row 41: The Write Stuff / Enchanted Mythril Ink
$ws.Range("H41").Value = 30304138
$ws.Range("I41").Value = 1196
$ws.Range("K41").Value = 1196
$ws.Range("M41").Value = -756

# ALC row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 7983.1665
$ws.Range("I62").Value = 7899
$ws.Range("K62").Value = 7899
$ws.Range("M62").Value = -7275

# ALC row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 7983.1665
$ws.Range("I65").Value = 7899
$ws.Range("K65").Value = 39495
$ws.Range("M65").Value = -36375

# ALC row 93: Spellbound / Koppranickel Index
$ws.Range("H93").Value = 46498.5
$ws.Range("J93").Value = 46498.5
$ws.Range("L93").Value = 46498.5
$ws.Range("N93").Value = -51490.5

# ALC row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 1515
$ws.Range("I98").Value = 1612.2778
$ws.Range("K98").Value = 1612.2778
$ws.Range("M98").Value = -114.2778000000001

# ALC row 99: Rumor Has It / Commanding Craftsman's Tea
$ws.Range("H99").Value = 27778062
$ws.Range("J99").Value = 417
$ws.Range("L99").Value = 1251
$ws.Range("N99").Value = -4247

# ALC row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 1515
$ws.Range("I122").Value = 1612.2778
$ws.Range("K122").Value = 4836.8334
$ws.Range("M122").Value = -2386.8334

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 7595.154
$ws.Range("I32").Value = 3918.5571
$ws.Range("K32").Value = 3918.5571
$ws.Range("M32").Value = -3631.5571

# ARM row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 7996275
$ws.Range("I45").Value = 14387403
$ws.Range("K45").Value = 14387403
$ws.Range("M45").Value = -14387026

# ARM row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 39634.957
$ws.Range("I74").Value = 25333.893
$ws.Range("J74").Value = 92548.89999999999
$ws.Range("K74").Value = 25333.893
$ws.Range("L74").Value = 92548.89999999999
$ws.Range("M74").Value = -24459.893
$ws.Range("N74").Value = -94296.89999999999

# ARM row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 39634.957
$ws.Range("I77").Value = 25333.893
$ws.Range("J77").Value = 92548.89999999999
$ws.Range("K77").Value = 126669.465
$ws.Range("L77").Value = 462744.5
$ws.Range("M77").Value = -122301.465
$ws.Range("N77").Value = -471480.5

# ARM row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1738571
$ws.Range("I102").Value = 2138925
$ws.Range("J102").Value = 3703.4443
$ws.Range("K102").Value = 2138925
$ws.Range("L102").Value = 3703.4443
$ws.Range("M102").Value = -2137303
$ws.Range("N102").Value = -6947.4443

# ARM row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 30388140
$ws.Range("J122").Value = 4175059.5
$ws.Range("L122").Value = 12525178.5
$ws.Range("N122").Value = -12530078.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 10754652
$ws.Range("J20").Value = 2578.182
$ws.Range("L20").Value = 2578.182
$ws.Range("N20").Value = -3072.182

# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 7151443
$ws.Range("I86").Value = 10011310
$ws.Range("J86").Value = 1774.75
$ws.Range("K86").Value = 10011310
$ws.Range("L86").Value = 1774.75
$ws.Range("M86").Value = -10010187
$ws.Range("N86").Value = -4020.75

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 7151443
$ws.Range("I89").Value = 10011310
$ws.Range("J89").Value = 1774.75
$ws.Range("K89").Value = 50056550
$ws.Range("L89").Value = 8873.75
$ws.Range("M89").Value = -50050934
$ws.Range("N89").Value = -20105.75

# BSM row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 20834064
$ws.Range("I105").Value = 20834064
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 20834064
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -20832317
$ws.Range("N105").ClearContents()

# BSM row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 2859639.8
$ws.Range("I107").Value = 3761246.5
$ws.Range("K107").Value = 3761246.5
$ws.Range("M107").Value = -3759326.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 2: In with the New / Bone Harpoon
$ws.Range("H2").Value = 9000
$ws.Range("I2").Value = 9000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -8887
$ws.Range("N2").ClearContents()

# CRP row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1619.6154
$ws.Range("I16").Value = 1289.625
$ws.Range("K16").Value = 1289.625
$ws.Range("M16").Value = -1002.625

# CRP row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1619.6154
$ws.Range("I113").Value = 1289.625
$ws.Range("K113").Value = 1289.625
$ws.Range("M113").Value = 880.375

# CRP row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 2546.1904
$ws.Range("I122").Value = 2188.5334
$ws.Range("J122").Value = 3440.3333
$ws.Range("K122").Value = 6565.600199999999
$ws.Range("L122").Value = 10320.9999
$ws.Range("M122").Value = -4115.600199999999
$ws.Range("N122").Value = -15220.9999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 34: Fever Pitch / Chamomile Tea
$ws.Range("H34").Value = 221.83333
$ws.Range("I34").Value = 221.83333
$ws.Range("K34").Value = 665.49999
$ws.Range("M34").Value = -581.49999

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 10420217
$ws.Range("J131").Value = 11498272
$ws.Range("L131").Value = 34494816
$ws.Range("N131").Value = -34504896

$ws = $wb.Worksheets.Item("GSM")
# GSM row 43: Get the Green Stuff / Malachite Earrings
$ws.Range("H43").Value = 2626
$ws.Range("I43").Value = 2151.2
$ws.Range("K43").Value = 2151.2
$ws.Range("M43").Value = -2000.2

# GSM row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 768874.2
$ws.Range("I97").Value = 1401247
$ws.Range("J97").Value = 992.9286
$ws.Range("K97").Value = 1401247
$ws.Range("L97").Value = 992.9286
$ws.Range("M97").Value = -1400751
$ws.Range("N97").Value = -1984.9286

# GSM row 99: Needle in a Hingan Stack / Dzo Horn Needle
$ws.Range("H99").Value = 3722.375
$ws.Range("I99").Value = 2825.5715
$ws.Range("K99").Value = 2825.5715
$ws.Range("M99").Value = -579.5715

# GSM row 104: Speak Softly and Carry a Metal Rod / Palladium Rod
$ws.Range("H104").Value = 76499.5
$ws.Range("J104").Value = 76499.5
$ws.Range("L104").Value = 76499.5
$ws.Range("N104").Value = -83487.5

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 525637
$ws.Range("I122").Value = 742993.75
$ws.Range("J122").Value = 3980.8
$ws.Range("K122").Value = 2228981.25
$ws.Range("L122").Value = 11942.4
$ws.Range("M122").Value = -2226531.25
$ws.Range("N122").Value = -16842.4

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 4157376.8
$ws.Range("J126").Value = 5211706.5
$ws.Range("L126").Value = 15635119.5
$ws.Range("N126").Value = -15640059.5

# GSM row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2423.6912
$ws.Range("I132").Value = 2312.4038
$ws.Range("J132").Value = 2785.375
$ws.Range("K132").Value = 6937.2114
$ws.Range("L132").Value = 8356.125
$ws.Range("M132").Value = -4407.2114
$ws.Range("N132").Value = -13416.125

$ws = $wb.Worksheets.Item("LTW")
# LTW row 29: Hands On / Fingerless Goatskin Gloves
$ws.Range("H29").Value = 24997.5
$ws.Range("I29").Value = 24997.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 24997.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -24702.5
$ws.Range("N29").ClearContents()

# LTW row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 5120.6772
$ws.Range("I100").Value = 5532.304
$ws.Range("J100").Value = 3937.25
$ws.Range("K100").Value = 5532.304
$ws.Range("L100").Value = 3937.25
$ws.Range("M100").Value = -4991.304
$ws.Range("N100").Value = -5019.25

$ws = $wb.Worksheets.Item("WVR")
# WVR row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 1823.7858
$ws.Range("J100").Value = 827.3333
$ws.Range("L100").Value = 1654.6666
$ws.Range("N100").Value = -2736.6666

# WVR row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 3466
$ws.Range("I122").Value = 2660
$ws.Range("J122").Value = 4675
$ws.Range("K122").Value = 7980
$ws.Range("L122").Value = 14025
$ws.Range("M122").Value = -5530
$ws.Range("N122").Value = -18925

# WVR row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1271.04
$ws.Range("I126").Value = 1266.4736
$ws.Range("J126").Value = 1285.5
$ws.Range("K126").Value = 3799.4208
$ws.Range("L126").Value = 3856.5
$ws.Range("M126").Value = -1329.4208
$ws.Range("N126").Value = -8796.5
